$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 ("Elon", 4) is being removed entirely - delete it so the used
# range / dimension shrinks back down instead of leaving a blank row.
$ws.Rows.Item(3).Delete()

# Headers
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "page_id"
$ws.Range("C1").Value = "emotion"
$ws.Range("D1").Value = "contact_details"

# Data row
$ws.Range("A2").Value = "b03e7319-7477-482f-96f2-cb89ddf0d08d"
$ws.Range("B2").Value = "pageid"
$ws.Range("C2").Value = ":-)"
